# Updates to column F ("想去人数" / "want-to-go count") across all four sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 147
$ws.Cells.Item(3, 6).Value = 966
$ws.Cells.Item(5, 6).Value = 2948
$ws.Cells.Item(6, 6).Value = 793
$ws.Cells.Item(8, 6).Value = 600
$ws.Cells.Item(9, 6).Value = 427
$ws.Cells.Item(11, 6).Value = 378
$ws.Cells.Item(12, 6).Value = 524
$ws.Cells.Item(13, 6).Value = 526
$ws.Cells.Item(14, 6).Value = 2157
$ws.Cells.Item(17, 6).Value = 18
$ws.Cells.Item(19, 6).Value = 2674
$ws.Cells.Item(25, 6).Value = 610
$ws.Cells.Item(28, 6).Value = 531
$ws.Cells.Item(29, 6).Value = 554
$ws.Cells.Item(30, 6).Value = 567
$ws.Cells.Item(31, 6).Value = 226
$ws.Cells.Item(32, 6).Value = 120
$ws.Cells.Item(33, 6).Value = 392
$ws.Cells.Item(34, 6).Value = 4675
$ws.Cells.Item(35, 6).Value = 242
$ws.Cells.Item(36, 6).Value = 25

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(8, 6).Value = 331
$ws.Cells.Item(9, 6).Value = 353
$ws.Cells.Item(14, 6).Value = 167
$ws.Cells.Item(23, 6).Value = 275
$ws.Cells.Item(24, 6).Value = 21
$ws.Cells.Item(25, 6).Value = 300
$ws.Cells.Item(27, 6).Value = 155
$ws.Cells.Item(31, 6).Value = 23
$ws.Cells.Item(36, 6).Value = 543

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(4, 6).Value = 1464
$ws.Cells.Item(5, 6).Value = 571
$ws.Cells.Item(6, 6).Value = 248
$ws.Cells.Item(7, 6).Value = 255

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 1464
$ws.Cells.Item(4, 6).Value = 571
$ws.Cells.Item(5, 6).Value = 147
$ws.Cells.Item(6, 6).Value = 248
$ws.Cells.Item(7, 6).Value = 966
$ws.Cells.Item(9, 6).Value = 2948
$ws.Cells.Item(10, 6).Value = 793
$ws.Cells.Item(12, 6).Value = 600
$ws.Cells.Item(13, 6).Value = 427
$ws.Cells.Item(16, 6).Value = 524
$ws.Cells.Item(17, 6).Value = 331
$ws.Cells.Item(18, 6).Value = 353
$ws.Cells.Item(19, 6).Value = 526
$ws.Cells.Item(21, 6).Value = 2157
$ws.Cells.Item(25, 6).Value = 167
$ws.Cells.Item(27, 6).Value = 2674
$ws.Cells.Item(35, 6).Value = 255
$ws.Cells.Item(37, 6).Value = 610
$ws.Cells.Item(38, 6).Value = 610
$ws.Cells.Item(39, 6).Value = 275
$ws.Cells.Item(40, 6).Value = 554
$ws.Cells.Item(41, 6).Value = 567
$ws.Cells.Item(42, 6).Value = 300
$ws.Cells.Item(43, 6).Value = 226
$ws.Cells.Item(45, 6).Value = 392
$ws.Cells.Item(47, 6).Value = 4675
$ws.Cells.Item(48, 6).Value = 242
$ws.Cells.Item(50, 6).Value = 543

